$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused rows 7-11 (old data rows no longer present)
$ws.Range("A7:C11").EntireRow.Delete()

# New header values for D1:F1, matching the style of B1:C1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("B1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)

# Updated data values for B2:F6
$data = @(
    @(44, 90, 8, 87, 16),
    @(95, 7, 60, 43, 27),
    @(72, 21, 55, 42, 45),
    @(12, 68, 28, 49, 39),
    @(59, 26, 25, 28, 51)
)

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
}
